# Append the 2025-05-03 price row to each of the daily price sheets in the
# Solar Prices workbook. Each sheet currently ends at row 62 (2025-05-02);
# we add row 63 with the new date and carry forward the prior day's price
# (matching the commit "Updated Argent prices in Excel").

$wb = $excel.ActiveWorkbook

$newDate = "2025-05-03"

$sheetUpdates = @(
    @{ Name = "N-Dense";                  Price = "38" },
    @{ Name = "N-Type";                   Price = "37.3" },
    @{ Name = "N-type Wafer";             Price = "1.02" },
    @{ Name = "Cell Topcon 183mm";        Price = "0.273" },
    @{ Name = "Module Topcon 183mm";      Price = "0.09" },
    @{ Name = "Silver Rear_side";         Price = "5,360" },
    @{ Name = "Silver Busbar front-side"; Price = "8,025" },
    @{ Name = "Silver finger front-side"; Price = "8,075" },
    @{ Name = "USD_CNY";                  Price = "7.2927" }
)

foreach ($update in $sheetUpdates) {
    $ws = $wb.Worksheets.Item($update.Name)
    $newRow = $ws.Cells.Item(62, 1).Row + 1

    $dateCell = $ws.Cells.Item($newRow, 1)
    # Prefix with an apostrophe so Excel stores the date as literal text
    # (matching the existing column A cells) instead of auto-converting
    # it to a date serial number, then reset the style so it doesn't pick
    # up the quote-prefix formatting.
    $dateCell.Value = "'" + $newDate
    $dateCell.Style = "Normal"

    $priceCell = $ws.Cells.Item($newRow, 2)
    # Same trick for the price, which is stored as text in this sheet too.
    $priceCell.Value = "'" + $update.Price
    $priceCell.Style = "Normal"
}
